# Added a pandas time series example to the data notebook: both sheets
# gain an extra "depth" row (4) with no temperature reading yet, pushing
# the final depth/temperature pair (5) down a row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet1 ("baseline") ---
# Insert a new row 5 (shifts old row 5: A5=5,B5=30 -> row 6), then fill
# the new row 5 with the depth-only reading.
$ws1.Rows.Item(5).Insert()
$ws1.Range("A5").Value = 4

# --- sheet2 ("perturbed") ---
$ws2.Rows.Item(5).Insert()
$ws2.Range("A5").Value = 4

# --- selections / active sheet ---
# Select sheet2's new last row first, then finish on sheet1's C24 so that
# sheet1 ends up the active tab (matches the workbook-level activeTab
# reverting to the default / first sheet).
[void]$ws2.Rows.Item(5).Select()
[void]$ws1.Range("C24").Select()
